$d = $word.ActiveDocument

# --- Update the date line ---
$d.Content.Find.Execute("2025-10-14 Tuesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-10-15 Wednesday", 2) | Out-Null

# --- Update the practice table ---
# The table has 20 rows x 5 columns; data lives in rows 1, 5, 9, 13, 17 (1-indexed).
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("11÷5=2, 1", "60÷5=12, 0", "31÷7=4, 3", "44÷6=7, 2", "70÷4=17, 2")
    5  = @("49÷5=9, 4", "22÷6=3, 4", "40÷8=5, 0", "45÷2=22, 1", "36÷3=12, 0")
    9  = @("88÷8=11, 0", "84÷6=14, 0", "73÷8=9, 1", "81÷9=9, 0", "51÷9=5, 6")
    13 = @("73÷2=36, 1", "80÷3=26, 2", "36÷6=6, 0", "81÷9=9, 0", "88÷7=12, 4")
    17 = @("60÷9=6, 6", "77÷6=12, 5", "36÷3=12, 0", "76÷2=38, 0", "39÷8=4, 7")
}

foreach ($rowIndex in $newValues.Keys) {
    $cols = $newValues[$rowIndex]
    for ($c = 1; $c -le $cols.Length; $c++) {
        $cell = $t.Cell($rowIndex, $c)
        $cell.Range.Text = $cols[$c - 1]
    }
}

Write-Output "edits applied"
